# Update marine fungi / protists data on the "Bochdansky" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bochdansky")

$ws.Range("B3").Value = 0.35
$ws.Range("C3").Value = 1.42

$ws.Range("B4").Value = 0.44
$ws.Range("C4").Value = 1.98

$ws.Range("B5").Value = 0.53
$ws.Range("C5").Value = 0.53

# Move the active selection to B5 (previously was C5)
$ws.Activate()
$ws.Range("B5").Select()
